$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.569.34"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "3.044.98"
$ws.Range("E3").Value = "  +4.20%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'200.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").Value = "'633.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.81%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.551"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +2.17%  "
$ws.Range("D10").Value = "3.043.82"
$ws.Range("E10").Value = "  +4.20%  "
$ws.Range("D11").Value = "'0.435"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.595.04"
$ws.Range("E13").Value = "  +3.99%  "
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").Value = "'5.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.23%  "
$ws.Range("D15").Value = "'29.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.68%  "
$ws.Range("D16").Value = "76.441.27"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "3.012.96"
$ws.Range("E18").Value = "  +3.32%  "
$ws.Range("D19").Value = "'13.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.27%  "
$ws.Range("D20").Value = "'9.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").Value = "'375.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("B22").Value = "Polkadot"
$ws.Range("C22").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D22").Value = "'4.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.14%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "'2.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("D24").Value = "'72.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.17%  "
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "'4.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.93%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "'9.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "'8.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.18%  "
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("D33").Value = "'512.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("E34").Value = "  +8.06%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'20.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.69%  "
$ws.Range("D37").Value = "'163.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").Value = "'0.387"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.67%  "
$ws.Range("D39").Value = "'20.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.52%  "
$ws.Range("B40").Value = "Cronos"
$ws.Range("C40").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D40").Value = "'0.105"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.47%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'187.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.97%  "
$ws.Range("E42").Value = "  -0.90%  "
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").Value = "'5.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").Value = "'43.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.25%  "
$ws.Range("D46").Value = "'1.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.88%  "
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").Value = "'0.615"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.60%  "
$ws.Range("D49").Value = "'0.717"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.15%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "'3.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.33%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'2.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.14%  "
